$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newParticipantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.instrument_model in ['DNBSEQ-G400']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id limit 100
'@

# Replace trailing newline introduced by the here-string with nothing; the text
# must not end with a newline character (matches original content exactly).
$newParticipantQuery = $newParticipantQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $newParticipantQuery

$ws.Range("C2").Select()
